$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-6: change date serial 45221 -> 45224
foreach ($row in 2..6) {
    $cell = $ws.Range("C$row")
    if ($cell.Value2 -eq 45221) {
        $cell.Value = 45224
    }
}
